{"js": "// Each entry is [oldText, newText] for one table-cell answer that the\n// commit regenerated (three-digit \u00f7 one-digit practice problems).\nconst replacements = [\n  [\"513\u00f76=85, 3\", \"984\u00f74=246, 0\"],\n  [\"845\u00f75=169, 0\", \"298\u00f74=74, 2\"],\n  [\"168\u00f78=21, 0\", \"368\u00f77=52, 4\"],\n  [\"122\u00f76=20, 2\", \"936\u00f74=234, 0\"],\n  [\"201\u00f73=67, 0\", \"830\u00f77=118, 4\"],\n  [\"543\u00f73=181, 0\", \"514\u00f75=102, 4\"],\n  [\"974\u00f77=139, 1\", \"401\u00f79=44, 5\"],\n  [\"410\u00f72=205, 0\", \"320\u00f77=45, 5\"],\n  [\"473\u00f78=59, 1\", \"964\u00f73=321, 1\"],\n  [\"929\u00f73=309, 2\", \"867\u00f78=108, 3\"],\n  [\"708\u00f77=101, 1\", \"518\u00f72=259, 0\"],\n  [\"790\u00f76=131, 4\", \"334\u00f74=83, 2\"],\n  [\"360\u00f76=60, 0\", \"931\u00f75=186, 1\"],\n  [\"855\u00f77=122, 1\", \"618\u00f72=309, 0\"],\n  [\"958\u00f79=106, 4\", \"589\u00f76=98, 1\"],\n  [\"143\u00f76=23, 5\", \"356\u00f76=59, 2\"],\n  [\"565\u00f74=141, 1\", \"810\u00f78=101, 2\"],\n  [\"678\u00f77=96, 6\", \"943\u00f72=471, 1\"],\n  [\"812\u00f75=162, 2\", \"992\u00f74=248, 0\"],\n  [\"193\u00f72=96, 1\", \"520\u00f79=57, 7\"],\n  [\"924\u00f76=154, 0\", \"164\u00f73=54, 2\"],\n  [\"888\u00f72=444, 0\", \"112\u00f79=12, 4\"],\n  [\"881\u00f72=440, 1\", \"115\u00f78=14, 3\"],\n  [\"875\u00f77=125, 0\", \"719\u00f72=359, 1\"],\n  [\"740\u00f74=185, 0\", \"550\u00f77=78, 4\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Each pair is the old/new answer text for one table-cell run that the\n# commit regenerated (three-digit \u00f7 one-digit practice problems).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"513\u00f76=85, 3\", \"984\u00f74=246, 0\"),\n    @(\"845\u00f75=169, 0\", \"298\u00f74=74, 2\"),\n    @(\"168\u00f78=21, 0\", \"368\u00f77=52, 4\"),\n    @(\"122\u00f76=20, 2\", \"936\u00f74=234, 0\"),\n    @(\"201\u00f73=67, 0\", \"830\u00f77=118, 4\"),\n    @(\"543\u00f73=181, 0\", \"514\u00f75=102, 4\"),\n    @(\"974\u00f77=139, 1\", \"401\u00f79=44, 5\"),\n    @(\"410\u00f72=205, 0\", \"320\u00f77=45, 5\"),\n    @(\"473\u00f78=59, 1\", \"964\u00f73=321, 1\"),\n    @(\"929\u00f73=309, 2\", \"867\u00f78=108, 3\"),\n    @(\"708\u00f77=101, 1\", \"518\u00f72=259, 0\"),\n    @(\"790\u00f76=131, 4\", \"334\u00f74=83, 2\"),\n    @(\"360\u00f76=60, 0\", \"931\u00f75=186, 1\"),\n    @(\"855\u00f77=122, 1\", \"618\u00f72=309, 0\"),\n    @(\"958\u00f79=106, 4\", \"589\u00f76=98, 1\"),\n    @(\"143\u00f76=23, 5\", \"356\u00f76=59, 2\"),\n    @(\"565\u00f74=141, 1\", \"810\u00f78=101, 2\"),\n    @(\"678\u00f77=96, 6\", \"943\u00f72=471, 1\"),\n    @(\"812\u00f75=162, 2\", \"992\u00f74=248, 0\"),\n    @(\"193\u00f72=96, 1\", \"520\u00f79=57, 7\"),\n    @(\"924\u00f76=154, 0\", \"164\u00f73=54, 2\"),\n    @(\"888\u00f72=444, 0\", \"112\u00f79=12, 4\"),\n    @(\"881\u00f72=440, 1\", \"115\u00f78=14, 3\"),\n    @(\"875\u00f77=125, 0\", \"719\u00f72=359, 1\"),\n    @(\"740\u00f74=185, 0\", \"550\u00f77=78, 4\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute(\n        $oldText,       # FindText\n        $false,         # MatchCase (text is distinctive enough without it)\n        $false,         # MatchWholeWord\n        $false,         # MatchWildcards\n        $false,         # MatchSoundsLike\n        $false,         # MatchAllWordForms\n        $true,          # Forward\n        1,              # Wrap (wdFindContinue)\n        $false,         # Format\n        $newText,       # ReplaceWith\n        2               # Replace (wdReplaceAll)\n    )\n    if (-not $found) {\n        throw \"No match found for: $oldText\"\n    }\n}\n"}
